$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 278209
$ws.Range("E10").Value = 1752512170
$ws.Range("C74").Value = 27992
$ws.Range("E74").Value = 54752531
$ws.Range("C100").Value = 9839
$ws.Range("E100").Value = 24692472
$ws.Range("C117").Value = 19731
$ws.Range("E117").Value = 56664962
$ws.Range("C168").Value = 285091
$ws.Range("E168").Value = 1212921927
$ws.Range("C169").Value = 562664
$ws.Range("E169").Value = 1286053236
$ws.Range("C170").Value = 367533
$ws.Range("E170").Value = 2847702736
$ws.Range("C171").Value = 115215
$ws.Range("E171").Value = 448609800
$ws.Range("C173").Value = 54396
$ws.Range("E173").Value = 151948108
$ws.Range("C174").Value = 357347
$ws.Range("E174").Value = 1019905253
$ws.Range("C175").Value = 125678
$ws.Range("E175").Value = 815152093
$ws.Range("C177").Value = 96778
$ws.Range("E177").Value = 174811932
$ws.Range("C179").Value = 235782
$ws.Range("E179").Value = 813286352
$ws.Range("C180").Value = 141523
$ws.Range("E180").Value = 341217178
$ws.Range("C188").Value = 19712
$ws.Range("E188").Value = 66204904
$ws.Range("C255").Value = 141370
$ws.Range("E255").Value = 414531469
$ws.Range("C280").Value = 95353
$ws.Range("E280").Value = 282552411
$ws.Range("C286").Value = 90609
$ws.Range("E286").Value = 162838526
$ws.Range("C313").Value = 220662
$ws.Range("E313").Value = 1371202181
$ws.Range("C322").Value = 81164
$ws.Range("E322").Value = 254557631
